# Adds a new "2021" data column (column R) to the table, mirroring the
# formatting of the existing column Q, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (thin header spacer row) - empty cell, same style as Q2
$ws.Range("Q2").Copy($ws.Range("R2"))

# Row 3 (year header row) - new year value 2021, same style as Q3
$ws.Range("Q3").Copy($ws.Range("R3"))
$ws.Range("R3").Value = 2021

# Row 4 (GVA share %, data row) - new value 13.5, same style as Q4
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 13.5

# Row 5 (GVA per capita, data row) - new value 15.1, same style as Q5
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 15.1

# Clear clipboard marching ants state left over from the Copy operations
$excel.CutCopyMode = $false

# Update the selected cell/range as recorded in the saved view
$ws.Range("T3").Select()
